# New crime data collected - weekly CompStat update (47th Precinct)
# Updates header metadata (Police Commissioner name, volume/issue number,
# reporting week dates) plus the weekly/28-day/YTD/2-year crime statistics
# table for rows 14-33.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header block
# ---------------------------------------------------------------------
$ws.Range("M6").Value = "Thomas G. Donlon"
$ws.Range("A8").Value = "Volume 31   Number  39"
$ws.Range("C9").Value = "Report Covering the Week  9/23/2024  Through  9/29/2024"

# ---------------------------------------------------------------------
# Row 14 - Murder
# ---------------------------------------------------------------------
$ws.Range("N14").Value = -78.260869565217

# ---------------------------------------------------------------------
# Row 15 - Rape
# ---------------------------------------------------------------------
$ws.Range("F15").Value = 4
$ws.Range("H15").Value = 300
$ws.Range("N15").Value = -43.859649122807

# ---------------------------------------------------------------------
# Row 16 - Robbery
# ---------------------------------------------------------------------
$ws.Range("C16").Value = 15
$ws.Range("D16").Value = 8
$ws.Range("E16").Value = 87.5
$ws.Range("F16").Value = 45
$ws.Range("H16").Value = 18.421052631578
$ws.Range("I16").Value = 397
$ws.Range("J16").Value = 352
$ws.Range("K16").Value = 12.784090909090
$ws.Range("L16").Value = 19.578313253012
$ws.Range("M16").Value = 28.478964401294
$ws.Range("N16").Value = -62.190476190476

# ---------------------------------------------------------------------
# Row 17 - Fel. Assault
# ---------------------------------------------------------------------
$ws.Range("C17").Value = 10
$ws.Range("D17").Value = 15
$ws.Range("E17").Value = -33.333333333333
$ws.Range("F17").Value = 59
$ws.Range("G17").Value = 79
$ws.Range("H17").Value = -25.316455696202
$ws.Range("I17").Value = 589
$ws.Range("J17").Value = 634
$ws.Range("K17").Value = -7.097791798107
$ws.Range("L17").Value = 5.366726296958
$ws.Range("M17").Value = 84.639498432601
$ws.Range("N17").Value = -12.351190476190

# ---------------------------------------------------------------------
# Row 18 - Burglary
# ---------------------------------------------------------------------
$ws.Range("C18").Value = 7
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = 133.333333333333
$ws.Range("F18").Value = 19
$ws.Range("G18").Value = 13
$ws.Range("H18").Value = 46.153846153846
$ws.Range("I18").Value = 206
$ws.Range("J18").Value = 206
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = -3.286384976525
$ws.Range("M18").Value = -17.6
$ws.Range("N18").Value = -85.472496473906

# ---------------------------------------------------------------------
# Row 19 - Gr. Larceny
# ---------------------------------------------------------------------
$ws.Range("C19").Value = 13
$ws.Range("E19").Value = -18.75
$ws.Range("F19").Value = 61
$ws.Range("G19").Value = 64
$ws.Range("H19").Value = -4.6875
$ws.Range("I19").Value = 675
$ws.Range("J19").Value = 580
$ws.Range("K19").Value = 16.379310344827
$ws.Range("L19").Value = 19.257950530035
$ws.Range("M19").Value = 175.510204081633
$ws.Range("N19").Value = 63.438256658595

# ---------------------------------------------------------------------
# Row 20 - G.L.A.
# ---------------------------------------------------------------------
$ws.Range("C20").Value = 10
$ws.Range("D20").Value = 13
$ws.Range("E20").Value = -23.076923076923
$ws.Range("F20").Value = 43
$ws.Range("G20").Value = 52
$ws.Range("H20").Value = -17.307692307692
$ws.Range("I20").Value = 410
$ws.Range("J20").Value = 483
$ws.Range("K20").Value = -15.113871635610
$ws.Range("L20").Value = 25
$ws.Range("M20").Value = 69.421487603305
$ws.Range("N20").Value = -66.255144032921

# ---------------------------------------------------------------------
# Row 21 - TOTAL
# ---------------------------------------------------------------------
$ws.Range("C21").Value = 55
$ws.Range("D21").Value = 55
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 232
$ws.Range("G21").Value = 247
$ws.Range("H21").Value = -6.072874493927
$ws.Range("I21").Value = 2314
$ws.Range("J21").Value = 2296
$ws.Range("K21").Value = 0.783972125435
$ws.Range("L21").Value = 12.98828125
$ws.Range("M21").Value = 63.533568904593
$ws.Range("N21").Value = -52.268976897689

# ---------------------------------------------------------------------
# Row 22 - Transit
# ---------------------------------------------------------------------
$ws.Range("F22").Value = 1
$ws.Range("H22").Value = -75
$ws.Range("J22").Value = 20
$ws.Range("K22").Value = 10
$ws.Range("L22").Value = -31.25
$ws.Range("M22").Value = 10

# ---------------------------------------------------------------------
# Row 23 - Housing
# ---------------------------------------------------------------------
$ws.Range("C23").Value = 3
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "***.*"
$ws.Range("G23").Value = 7
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 83
$ws.Range("K23").Value = 3.75
$ws.Range("L23").Value = -6.741573033707
$ws.Range("M23").Value = 43.103448275862

# ---------------------------------------------------------------------
# Row 24 - Petit Larceny
# ---------------------------------------------------------------------
$ws.Range("C24").Value = 22
$ws.Range("D24").Value = 22
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 98
$ws.Range("G24").Value = 93
$ws.Range("H24").Value = 5.376344086021
$ws.Range("I24").Value = 967
$ws.Range("J24").Value = 1005
$ws.Range("K24").Value = -3.781094527363
$ws.Range("L24").Value = -14.348981399468
$ws.Range("M24").Value = 71.453900709219

# ---------------------------------------------------------------------
# Row 25 - Retail Theft
# ---------------------------------------------------------------------
$ws.Range("C25").Value = 10
$ws.Range("D25").Value = 2
$ws.Range("E25").Value = 400
$ws.Range("F25").Value = 42
$ws.Range("G25").Value = 20
$ws.Range("H25").Value = 110
$ws.Range("I25").Value = 296
$ws.Range("J25").Value = 276
$ws.Range("K25").Value = 7.246376811594
$ws.Range("L25").Value = -13.702623906705

# ---------------------------------------------------------------------
# Row 26 - Misd. Assault
# ---------------------------------------------------------------------
$ws.Range("D26").Value = 16
$ws.Range("E26").Value = 18.75
$ws.Range("F26").Value = 78
$ws.Range("G26").Value = 67
$ws.Range("H26").Value = 16.417910447761
$ws.Range("I26").Value = 845
$ws.Range("J26").Value = 712
$ws.Range("K26").Value = 18.679775280898
$ws.Range("L26").Value = 24.631268436578
$ws.Range("M26").Value = 13.422818791946

# ---------------------------------------------------------------------
# Row 27 - UCR Rape*
# ---------------------------------------------------------------------
$ws.Range("D27").Value = 3
$ws.Range("F27").Value = 4
$ws.Range("G27").Value = 7
$ws.Range("H27").Value = -42.857142857142
$ws.Range("J27").Value = 50
$ws.Range("K27").Value = -6

# ---------------------------------------------------------------------
# Row 28 - Other Sex Crimes
# ---------------------------------------------------------------------
$ws.Range("C28").NumberFormat = "#,##0"
$ws.Range("C28").Value = 1
$ws.Range("D28").Value = 4
$ws.Range("E28").Value = -75
$ws.Range("F28").Value = 8
$ws.Range("H28").Value = -11.111111111111
$ws.Range("I28").Value = 66
$ws.Range("J28").Value = 62
$ws.Range("K28").Value = 6.451612903225
$ws.Range("L28").Value = 11.864406779661

# ---------------------------------------------------------------------
# Row 29 - Shooting Vic.
# ---------------------------------------------------------------------
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "***.*"
$ws.Range("L29").Value = -35.135135135135
$ws.Range("M29").Value = -54.716981132075
$ws.Range("N29").Value = -78.947368421052

# ---------------------------------------------------------------------
# Row 30 - Shooting Inc.
# ---------------------------------------------------------------------
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "***.*"
$ws.Range("L30").Value = -52.941176470588
$ws.Range("M30").Value = -62.790697674418
$ws.Range("N30").Value = -85.046728971962

# ---------------------------------------------------------------------
# Row 33 - Traffic Fatalities
# ---------------------------------------------------------------------
$ws.Range("C33").NumberFormat = "#,##0"
$ws.Range("C33").Value = 1
$ws.Range("F33").NumberFormat = "#,##0"
$ws.Range("F33").Value = 1
$ws.Range("I33").NumberFormat = "#,##0"
$ws.Range("I33").Value = 1
$ws.Range("K33").Value = -50
$ws.Range("L33").Value = -50
